# Updates cryptos list values (price + 1h volume/change columns) to match
# the latest scrape, per commit "Updated cryptos list on Sun May 26
# 15:48:10 UTC 2024 with GitHub Actions".
#
# Column D ("Price") and E ("Volume(1h)") cells are stored as plain text in
# the sheet (prices use "." as a thousands separator, e.g. "69.149.99", and
# percentages are padded with double spaces, e.g. "  +0.07%  "). Whenever the
# new value is a bare decimal number (no thousands separator) we briefly force
# the cell to Text format so Excel's input parser does not reinterpret it as
# a numeric value, then restore the cell's original (Normal) style so no
# formatting residue is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.149.99"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "3.858.97"
$ws.Range("E3").Value = "  +2.91%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.25%  "
$ws.Range("D7").Value = "3.858.67"
$ws.Range("E7").Value = "  +2.99%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.36%  "
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.34"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.14"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000246"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").Value = "4.514.11"
$ws.Range("E15").Value = "  +3.10%  "
$ws.Range("D16").Value = "3.864.31"
$ws.Range("E16").Value = "  +3.24%  "
$ws.Range("D17").Value = "69.341.87"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.90%  "
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "489.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.724"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000159"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.29%  "
$ws.Range("E26").Value = "  -1.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("E29").Value = "  -1.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.33%  "
$ws.Range("D33").Value = "4.023.43"
$ws.Range("E33").Value = "  +3.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "32.41"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("D35").Value = "3.814.06"
$ws.Range("E35").Value = "  +3.49%  "
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("E37").Value = "  +2.42%  "
$ws.Range("E38").Value = "  +4.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.320"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "441.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.05%  "
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("B46").Value = "Cosmos"
$ws.Range("C46").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.07%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +18.65%  "
$ws.Range("D49").Value = "2.860.73"
$ws.Range("E49").Value = "  +2.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "143.11"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0359"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.79%  "
